$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.929.01"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "1.665.32"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.65"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.532"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.65%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.27"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0897"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.92%  "
$ws.Range("D12").Value = "1.900.16"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "1.648.60"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.07"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.14"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "26.928.46"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "234.08"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.02"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.37"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.12"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.15"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.89"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("D33").Value = "1.455.29"
$ws.Range("E33").Value = "  -4.85%  "
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.63"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.67%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.580"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.73"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.52%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.974"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.90"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("D45").Value = "1.810.98"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.783"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.39"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "0.0₆0104"
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.101"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.81%  "
$ws.Range("E51").Value = "  +0.26%  "
